# New .ttl from Google sheet has been generated
# - D26:D78 lose the ",skos:Concept" suffix (now just "iop:VariableSet")
# - U26's modified-date moves from 2023-08-23 to 2023-09-27
# - column AF (all-empty placeholder column) is removed entirely
# - row 84 (a stray extra row) is removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 26-78: "iop:VariableSet,skos:Concept" -> "iop:VariableSet"
$ws.Range("D26:D78").Value = "iop:VariableSet"

# U26: 2023-08-23 -> 2023-09-27 (force text so it isn't coerced into a date)
$dateCell = $ws.Cells.Item(26, 21)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2023-09-27"
$dateCell.Style = "Normal"

# Drop the now-empty AF column entirely
$ws.Columns("AF").Delete()

# Drop the trailing row 84
$ws.Rows(84).Delete()
